$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9 (ICSA_thou - Initial Jobless Claims) data refresh
$ws.Range("E9").Value = 212000
$ws.Range("G9").Value = 363670.4980842912
$ws.Range("H9").Value = -12000
$ws.Range("I9").Value = -0.05357142857142857
